$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells: Q1 "Purchase Count", R1 "Multiple Purchases" ---
$ws.Range("Q1").Value = "Purchase Count"
$ws.Range("R1").Value = "Multiple Purchases"

# Match the header formatting (bold / centered / thin border) used by the
# other header cells (e.g. P1) by copying its format onto the new headers.
$ws.Range("P1").Copy()
$ws.Range("Q1:R1").PasteSpecial(-4122)

# --- New data columns per row ---
# Q = "Purchase Count": number of purchases tied to the customer's email.
# Every email in this sheet is unique, so the count is 1 for every row that
# has an email address; rows with a missing email (5 and 12) are left blank.
# R = "Multiple Purchases": TRUE when Purchase Count > 1 - FALSE everywhere
# here since there are no repeated emails.

$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = $false

$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = $false

$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = $false

# Row 5: Email column is blank, so Purchase Count is left blank too.
$ws.Range("R5").Value = $false

$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = $false

$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = $false

$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = $false

$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = $false

$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = $false

$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = $false

# Row 12: Email column is blank, so Purchase Count is left blank too.
$ws.Range("R12").Value = $false

$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = $false

$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = $false

$ws.Range("Q15").Value = 1
$ws.Range("R15").Value = $false
